# Update countries & provincias Spain
#
# This applies the daily data refresh to the "Pais" sheet:
#  - Updated "Datos actualizados..." timestamp
#  - Updated totals for Estados Unidos (row 4) and Burkina Faso (row 88)
#  - Australia's row moved up (now before Brasil) with refreshed Australia
#    figures, Brasil's figures carried down unchanged, and Suecia's figures
#    shifted down with a small update
#  - Trinidad yTobago's row moved up (now before Nigeria) with refreshed
#    Trinidad yTobago figures, Nigeria's figures carried down unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 22:42"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 82179
$ws.Range("C4").Value = 13968
$ws.Range("E4").Value = 79138

# Reorder: Australia now sits above Brasil, Suecia drops below Brasil
# Row 21 -> Australia (new figures)
$ws.Range("A21").Value = "Australia"
$ws.Range("B21").Value = 2985
$ws.Range("C21").Value = 309
$ws.Range("D21").Value = 170
$ws.Range("E21").Value = 2802
$ws.Range("F21").Value = 11
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 13

# Row 22 -> Brasil (figures carried down unchanged from old row 21)
$ws.Range("A22").Value = "Brasil"
$ws.Range("B22").Value = 2915
$ws.Range("C22").Value = 361
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 2832
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = 77

# Row 23 -> Suecia (mostly carried down, with a small refresh)
$ws.Range("A23").Value = "Suecia"
$ws.Range("B23").Value = 2840
$ws.Range("C23").Value = 314
$ws.Range("D23").Value = 16
$ws.Range("E23").Value = 2747
$ws.Range("F23").Value = 176
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 77

# Burkina Faso (row 88)
$ws.Range("E88").Value = 135
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 7

# Reorder: Trinidad yTobago now sits above Nigeria
# Row 114 -> Trinidad yTobago (new figures)
$ws.Range("A114").Value = "Trinidad yTobago"
$ws.Range("B114").Value = 65
$ws.Range("C114").Value = 5
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 64

# Row 115 -> Nigeria (figures carried down unchanged from old row 114)
$ws.Range("A115").Value = "Nigeria"
$ws.Range("B115").Value = 65
$ws.Range("C115").Value = 14
$ws.Range("D115").Value = 2
$ws.Range("E115").Value = 62
